$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.394.84"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "3.379.43"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'574.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").Value = "'137.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.60%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.378.48"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").Value = "'0.473"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("D10").Value = "'7.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.72%  "

$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").Value = "'0.388"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("D13").Value = "3.955.71"
$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("E14").Value = "  +2.38%  "

$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").Value = "'26.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.60%  "

$ws.Range("D17").Value = "3.380.71"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("D18").Value = "61.538.15"
$ws.Range("E18").Value = "  +0.09%  "

$ws.Range("D19").Value = "'14.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("D20").Value = "'5.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("E21").Value = "  -1.23%  "

$ws.Range("D22").Value = "'376.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.58%  "

$ws.Range("E23").Value = "  -3.76%  "

$ws.Range("D24").Value = "3.512.20"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").Value = "'0.0000126"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.86%  "

$ws.Range("D27").Value = "'71.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("D28").Value = "'1.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.41%  "

$ws.Range("E29").Value = "  -3.95%  "

$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").Value = "'8.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.92%  "

$ws.Range("E32").Value = "  +2.72%  "

$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").Value = "'23.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").Value = "'5.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.52%  "

$ws.Range("D37").Value = "'6.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.55%  "

$ws.Range("E38").Value = "  -1.29%  "

$ws.Range("D39").Value = "'166.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.90%  "

$ws.Range("D40").Value = "'0.0772"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.25%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("E42").Value = "  +2.29%  "

$ws.Range("E43").Value = "  +1.64%  "

$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("E45").Value = "  -0.79%  "

$ws.Range("D46").Value = "'1.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("D47").Value = "'24.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.41%  "

$ws.Range("D48").Value = "'6.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.86%  "

$ws.Range("D49").Value = "'22.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.68%  "

$ws.Range("D50").Value = "2.364.48"
$ws.Range("E50").Value = "  +3.80%  "

$ws.Range("E51").Value = "  -1.31%  "
